$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header E1
$ws.Range("E1").Value = "Choices"

# Table data (rows 2-11): Sno, Name, Roll Number, Email id, Choices
$colA = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$colB = @("ASDF", "qdy dfef", "qrtyw ", "aassw", "weff_", "eef", "deed", "dffc", "eefv", "efccc")
$colC = @(45, 46, 47, 48, 49, 50, 51, 52, 53, 54)
$colD = @("qwertyuio", "asdfghjkl", "sdfghjkl;", "dfghjkl", "poiuytgfrd", "xcvbnm,", "dfghjkl", "xcvbnm,.", "xcvbnm,.", "plkjnhgfd")
$colE = @("A,B,C", "A", "B", "C", "D", "E", "A,B,C,D", "A,B,C,D,E", "E", "B, C")

# Shared-string indices are allocated in the order cells are populated,
# reproducing the original authoring order: column A, then column B (top
# to bottom), then column C, then column D (top to bottom), then column E
# -- but column E values were entered in a non-sequential row order.
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

$eRowOrder = @(2, 3, 6, 4, 5, 8, 7, 9, 10, 11)
foreach ($r in $eRowOrder) {
    $ws.Cells.Item($r, 5).Value = $colE[$r - 2]
}

# Column widths (target widths 19.6640625 / 19.109375 chars as stored in
# the XML "width" attribute, which already includes cell padding)
$ws.Range("C1").ColumnWidth = 18.833333333333332
$ws.Range("D1").ColumnWidth = 18.333333333333332

# Selection
$ws.Range("I16").Select()
